$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 7 through 29 (content beyond the new A1:A6 range)
$ws.Range("A7:A29").EntireRow.Delete()

# Update rows 2-6 with the combined card text representations
$ws.Range("A2").Value = "('Killing Wave', ['{X}{B}', 'Sorcery', 'For each creature, its controller sacrifices it unless they pay X life.'])"
$ws.Range("A3").Value = "('Latch Seeker', ['{1}{U}{U}', 'Creature — Spirit', 'Latch Seeker can" + [char]8217 + "t be blocked.', '3/1'])"
$ws.Range("A4").Value = "('Moonsilver Spear', ['{4}', 'Artifact — Equipment', 'Equipped creature has first strike.', 'Whenever equipped creature attacks, create a 4/4 white Angel creature token with flying.', 'Equip {4}'])"
$ws.Range("A5").Value = "('Restoration Angel', ['{3}{W}', 'Creature — Angel', 'Flash', 'Flying', 'When Restoration Angel enters the battlefield, you may exile target non-Angel creature you control, then return that card to the battlefield under your control.', '3/4'])"
$ws.Range("A6").Value = "('Silverblade Paladin', ['{1}{W}{W}', 'Creature — Human Knight', 'Soulbond (You may pair this creature with another unpaired creature when either enters the battlefield. They remain paired for as long as you control both of them.)', 'As long as Silverblade Paladin is paired with another creature, both creatures have double strike.', '2/2'])"
